$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append 13 new row(s) ----
$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = '2026-01-28'
$ws.Cells.Item(40, 1).Style = "Normal"
$ws.Cells.Item(40, 2).Value = '14:42:46'
$ws.Cells.Item(40, 3).Value = '14:00'
$ws.Cells.Item(40, 4).Value = 'Bathroom'
$ws.Cells.Item(40, 5).Value = 'No Motion'
$ws.Cells.Item(40, 6).Value = 'Inactive'
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = '2026-01-28'
$ws.Cells.Item(41, 1).Style = "Normal"
$ws.Cells.Item(41, 2).Value = '14:42:48'
$ws.Cells.Item(41, 3).Value = '14:00'
$ws.Cells.Item(41, 4).Value = 'Bathroom'
$ws.Cells.Item(41, 5).Value = 'No Motion'
$ws.Cells.Item(41, 6).Value = 'Inactive'
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = '2026-01-28'
$ws.Cells.Item(42, 1).Style = "Normal"
$ws.Cells.Item(42, 2).Value = '14:42:53'
$ws.Cells.Item(42, 3).Value = '14:00'
$ws.Cells.Item(42, 4).Value = 'Bathroom'
$ws.Cells.Item(42, 5).Value = 'No Motion'
$ws.Cells.Item(42, 6).Value = 'Inactive'
$ws.Cells.Item(43, 1).NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = '2026-01-28'
$ws.Cells.Item(43, 1).Style = "Normal"
$ws.Cells.Item(43, 2).Value = '14:42:58'
$ws.Cells.Item(43, 3).Value = '14:00'
$ws.Cells.Item(43, 4).Value = 'Bathroom'
$ws.Cells.Item(43, 5).Value = 'No Motion'
$ws.Cells.Item(43, 6).Value = 'Inactive'
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = '2026-01-28'
$ws.Cells.Item(44, 1).Style = "Normal"
$ws.Cells.Item(44, 2).Value = '14:43:03'
$ws.Cells.Item(44, 3).Value = '14:00'
$ws.Cells.Item(44, 4).Value = 'Bathroom'
$ws.Cells.Item(44, 5).Value = 'No Motion'
$ws.Cells.Item(44, 6).Value = 'Inactive'
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = '2026-01-28'
$ws.Cells.Item(45, 1).Style = "Normal"
$ws.Cells.Item(45, 2).Value = '14:43:08'
$ws.Cells.Item(45, 3).Value = '14:00'
$ws.Cells.Item(45, 4).Value = 'Bathroom'
$ws.Cells.Item(45, 5).Value = 'No Motion'
$ws.Cells.Item(45, 6).Value = 'Inactive'
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = '2026-01-28'
$ws.Cells.Item(46, 1).Style = "Normal"
$ws.Cells.Item(46, 2).Value = '14:43:13'
$ws.Cells.Item(46, 3).Value = '14:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).Value = 'No Motion'
$ws.Cells.Item(46, 6).Value = 'Inactive'
$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = '2026-01-28'
$ws.Cells.Item(47, 1).Style = "Normal"
$ws.Cells.Item(47, 2).Value = '14:43:18'
$ws.Cells.Item(47, 3).Value = '14:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).Value = 'No Motion'
$ws.Cells.Item(47, 6).Value = 'Inactive'
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = '2026-01-28'
$ws.Cells.Item(48, 1).Style = "Normal"
$ws.Cells.Item(48, 2).Value = '14:43:23'
$ws.Cells.Item(48, 3).Value = '14:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).Value = 'No Motion'
$ws.Cells.Item(48, 6).Value = 'Inactive'
$ws.Cells.Item(49, 1).NumberFormat = "@"
$ws.Cells.Item(49, 1).Value = '2026-01-28'
$ws.Cells.Item(49, 1).Style = "Normal"
$ws.Cells.Item(49, 2).Value = '14:43:28'
$ws.Cells.Item(49, 3).Value = '14:00'
$ws.Cells.Item(49, 4).Value = 'Bathroom'
$ws.Cells.Item(49, 5).Value = 'No Motion'
$ws.Cells.Item(49, 6).Value = 'Inactive'
$ws.Cells.Item(50, 1).NumberFormat = "@"
$ws.Cells.Item(50, 1).Value = '2026-01-28'
$ws.Cells.Item(50, 1).Style = "Normal"
$ws.Cells.Item(50, 2).Value = '14:43:33'
$ws.Cells.Item(50, 3).Value = '14:00'
$ws.Cells.Item(50, 4).Value = 'Bathroom'
$ws.Cells.Item(50, 5).Value = 'No Motion'
$ws.Cells.Item(50, 6).Value = 'Inactive'
$ws.Cells.Item(51, 1).NumberFormat = "@"
$ws.Cells.Item(51, 1).Value = '2026-01-28'
$ws.Cells.Item(51, 1).Style = "Normal"
$ws.Cells.Item(51, 2).Value = '14:43:38'
$ws.Cells.Item(51, 3).Value = '14:00'
$ws.Cells.Item(51, 4).Value = 'Bathroom'
$ws.Cells.Item(51, 5).Value = 'No Motion'
$ws.Cells.Item(51, 6).Value = 'Inactive'
$ws.Cells.Item(52, 1).NumberFormat = "@"
$ws.Cells.Item(52, 1).Value = '2026-01-28'
$ws.Cells.Item(52, 1).Style = "Normal"
$ws.Cells.Item(52, 2).Value = '14:43:43'
$ws.Cells.Item(52, 3).Value = '14:00'
$ws.Cells.Item(52, 4).Value = 'Bathroom'
$ws.Cells.Item(52, 5).Value = 'No Motion'
$ws.Cells.Item(52, 6).Value = 'Inactive'

# ---- Humidity sheet: append 12 new row(s) ----
$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = '2026-01-28'
$ws.Cells.Item(37, 1).Style = "Normal"
$ws.Cells.Item(37, 2).Value = '14:42:46'
$ws.Cells.Item(37, 3).Value = '14:00'
$ws.Cells.Item(37, 4).Value = 'Bathroom'
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '88.7%'
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(37, 6).Value = 'Active'
$ws.Cells.Item(38, 1).NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = '2026-01-28'
$ws.Cells.Item(38, 1).Style = "Normal"
$ws.Cells.Item(38, 2).Value = '14:42:48'
$ws.Cells.Item(38, 3).Value = '14:00'
$ws.Cells.Item(38, 4).Value = 'Bathroom'
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '88.7%'
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(38, 6).Value = 'Active'
$ws.Cells.Item(39, 1).NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = '2026-01-28'
$ws.Cells.Item(39, 1).Style = "Normal"
$ws.Cells.Item(39, 2).Value = '14:42:52'
$ws.Cells.Item(39, 3).Value = '14:00'
$ws.Cells.Item(39, 4).Value = 'Bathroom'
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '87.8%'
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(39, 6).Value = 'Active'
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = '2026-01-28'
$ws.Cells.Item(40, 1).Style = "Normal"
$ws.Cells.Item(40, 2).Value = '14:43:00'
$ws.Cells.Item(40, 3).Value = '14:00'
$ws.Cells.Item(40, 4).Value = 'Bathroom'
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '87.8%'
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(40, 6).Value = 'Active'
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = '2026-01-28'
$ws.Cells.Item(41, 1).Style = "Normal"
$ws.Cells.Item(41, 2).Value = '14:43:04'
$ws.Cells.Item(41, 3).Value = '14:00'
$ws.Cells.Item(41, 4).Value = 'Bathroom'
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '88.7%'
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(41, 6).Value = 'Active'
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = '2026-01-28'
$ws.Cells.Item(42, 1).Style = "Normal"
$ws.Cells.Item(42, 2).Value = '14:43:12'
$ws.Cells.Item(42, 3).Value = '14:00'
$ws.Cells.Item(42, 4).Value = 'Bathroom'
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '87.7%'
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(42, 6).Value = 'Active'
$ws.Cells.Item(43, 1).NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = '2026-01-28'
$ws.Cells.Item(43, 1).Style = "Normal"
$ws.Cells.Item(43, 2).Value = '14:43:16'
$ws.Cells.Item(43, 3).Value = '14:00'
$ws.Cells.Item(43, 4).Value = 'Bathroom'
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '88.7%'
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 6).Value = 'Active'
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = '2026-01-28'
$ws.Cells.Item(44, 1).Style = "Normal"
$ws.Cells.Item(44, 2).Value = '14:43:24'
$ws.Cells.Item(44, 3).Value = '14:00'
$ws.Cells.Item(44, 4).Value = 'Bathroom'
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '87.8%'
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 6).Value = 'Active'
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = '2026-01-28'
$ws.Cells.Item(45, 1).Style = "Normal"
$ws.Cells.Item(45, 2).Value = '14:43:28'
$ws.Cells.Item(45, 3).Value = '14:00'
$ws.Cells.Item(45, 4).Value = 'Bathroom'
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '88.7%'
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(45, 6).Value = 'Active'
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = '2026-01-28'
$ws.Cells.Item(46, 1).Style = "Normal"
$ws.Cells.Item(46, 2).Value = '14:43:32'
$ws.Cells.Item(46, 3).Value = '14:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '87.8%'
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(46, 6).Value = 'Active'
$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = '2026-01-28'
$ws.Cells.Item(47, 1).Style = "Normal"
$ws.Cells.Item(47, 2).Value = '14:43:36'
$ws.Cells.Item(47, 3).Value = '14:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '88.7%'
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(47, 6).Value = 'Active'
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = '2026-01-28'
$ws.Cells.Item(48, 1).Style = "Normal"
$ws.Cells.Item(48, 2).Value = '14:43:44'
$ws.Cells.Item(48, 3).Value = '14:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '88.7%'
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(48, 6).Value = 'Active'

# ---- Temperature sheet: append 13 new row(s) ----
$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(36, 1).NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = '2026-01-28'
$ws.Cells.Item(36, 1).Style = "Normal"
$ws.Cells.Item(36, 2).Value = '14:42:45'
$ws.Cells.Item(36, 3).Value = '14:00'
$ws.Cells.Item(36, 4).Value = 'Bathroom'
$ws.Cells.Item(36, 5).Value = '22.7C'
$ws.Cells.Item(36, 6).Value = 'Active'
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = '2026-01-28'
$ws.Cells.Item(37, 1).Style = "Normal"
$ws.Cells.Item(37, 2).Value = '14:42:46'
$ws.Cells.Item(37, 3).Value = '14:00'
$ws.Cells.Item(37, 4).Value = 'Bathroom'
$ws.Cells.Item(37, 5).Value = '22.6C'
$ws.Cells.Item(37, 6).Value = 'Active'
$ws.Cells.Item(38, 1).NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = '2026-01-28'
$ws.Cells.Item(38, 1).Style = "Normal"
$ws.Cells.Item(38, 2).Value = '14:42:48'
$ws.Cells.Item(38, 3).Value = '14:00'
$ws.Cells.Item(38, 4).Value = 'Bathroom'
$ws.Cells.Item(38, 5).Value = '22.7C'
$ws.Cells.Item(38, 6).Value = 'Active'
$ws.Cells.Item(39, 1).NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = '2026-01-28'
$ws.Cells.Item(39, 1).Style = "Normal"
$ws.Cells.Item(39, 2).Value = '14:42:52'
$ws.Cells.Item(39, 3).Value = '14:00'
$ws.Cells.Item(39, 4).Value = 'Bathroom'
$ws.Cells.Item(39, 5).Value = '22.7C'
$ws.Cells.Item(39, 6).Value = 'Active'
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = '2026-01-28'
$ws.Cells.Item(40, 1).Style = "Normal"
$ws.Cells.Item(40, 2).Value = '14:43:00'
$ws.Cells.Item(40, 3).Value = '14:00'
$ws.Cells.Item(40, 4).Value = 'Bathroom'
$ws.Cells.Item(40, 5).Value = '22.7C'
$ws.Cells.Item(40, 6).Value = 'Active'
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = '2026-01-28'
$ws.Cells.Item(41, 1).Style = "Normal"
$ws.Cells.Item(41, 2).Value = '14:43:04'
$ws.Cells.Item(41, 3).Value = '14:00'
$ws.Cells.Item(41, 4).Value = 'Bathroom'
$ws.Cells.Item(41, 5).Value = '22.7C'
$ws.Cells.Item(41, 6).Value = 'Active'
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = '2026-01-28'
$ws.Cells.Item(42, 1).Style = "Normal"
$ws.Cells.Item(42, 2).Value = '14:43:12'
$ws.Cells.Item(42, 3).Value = '14:00'
$ws.Cells.Item(42, 4).Value = 'Bathroom'
$ws.Cells.Item(42, 5).Value = '22.7C'
$ws.Cells.Item(42, 6).Value = 'Active'
$ws.Cells.Item(43, 1).NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = '2026-01-28'
$ws.Cells.Item(43, 1).Style = "Normal"
$ws.Cells.Item(43, 2).Value = '14:43:16'
$ws.Cells.Item(43, 3).Value = '14:00'
$ws.Cells.Item(43, 4).Value = 'Bathroom'
$ws.Cells.Item(43, 5).Value = '22.7C'
$ws.Cells.Item(43, 6).Value = 'Active'
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = '2026-01-28'
$ws.Cells.Item(44, 1).Style = "Normal"
$ws.Cells.Item(44, 2).Value = '14:43:24'
$ws.Cells.Item(44, 3).Value = '14:00'
$ws.Cells.Item(44, 4).Value = 'Bathroom'
$ws.Cells.Item(44, 5).Value = '22.7C'
$ws.Cells.Item(44, 6).Value = 'Active'
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = '2026-01-28'
$ws.Cells.Item(45, 1).Style = "Normal"
$ws.Cells.Item(45, 2).Value = '14:43:28'
$ws.Cells.Item(45, 3).Value = '14:00'
$ws.Cells.Item(45, 4).Value = 'Bathroom'
$ws.Cells.Item(45, 5).Value = '22.7C'
$ws.Cells.Item(45, 6).Value = 'Active'
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = '2026-01-28'
$ws.Cells.Item(46, 1).Style = "Normal"
$ws.Cells.Item(46, 2).Value = '14:43:32'
$ws.Cells.Item(46, 3).Value = '14:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).Value = '22.7C'
$ws.Cells.Item(46, 6).Value = 'Active'
$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = '2026-01-28'
$ws.Cells.Item(47, 1).Style = "Normal"
$ws.Cells.Item(47, 2).Value = '14:43:36'
$ws.Cells.Item(47, 3).Value = '14:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).Value = '22.7C'
$ws.Cells.Item(47, 6).Value = 'Active'
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = '2026-01-28'
$ws.Cells.Item(48, 1).Style = "Normal"
$ws.Cells.Item(48, 2).Value = '14:43:44'
$ws.Cells.Item(48, 3).Value = '14:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).Value = '22.7C'
$ws.Cells.Item(48, 6).Value = 'Active'

# ---- Proximity sheet: append 2 new row(s) ----
$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = '2026-01-28'
$ws.Cells.Item(17, 1).Style = "Normal"
$ws.Cells.Item(17, 2).Value = '14:42:45'
$ws.Cells.Item(17, 3).Value = '14:00'
$ws.Cells.Item(17, 4).Value = 'Bathroom Door'
$ws.Cells.Item(17, 5).Value = 'EXIT'
$ws.Cells.Item(17, 6).Value = 'User EXITED Bathroom'
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = '2026-01-28'
$ws.Cells.Item(18, 1).Style = "Normal"
$ws.Cells.Item(18, 2).Value = '14:42:46'
$ws.Cells.Item(18, 3).Value = '14:00'
$ws.Cells.Item(18, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(18, 5).Value = 'EXIT'
$ws.Cells.Item(18, 6).Value = 'User EXITED Living Room Main Door'

# ---- Camera sheet: append 1 new row(s) ----
$ws = $wb.Worksheets.Item("Camera")
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = '2026-01-28'
$ws.Cells.Item(9, 1).Style = "Normal"
$ws.Cells.Item(9, 2).Value = '14:42:45'
$ws.Cells.Item(9, 3).Value = '14:00'
$ws.Cells.Item(9, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(9, 5).Value = 'Image Captured'
$ws.Cells.Item(9, 6).Value = 'Active'
